$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "pooja168@givmail.com"
$ws.Range("D7").Select()
